$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Premium" / "10 cm" data row (currently row 14) in place ---
# Fix the thickness label formatting ("10 cm" -> "10cm") and the package/total counts.
$ws.Range("B14").Value = "10cm"
$ws.Range("C14").Value = 100
$ws.Range("D14").Value = 1095

# --- Remove the remaining product rows (old rows 15:23) ---
$ws.Rows("15:23").Delete()

# --- Remove the two blank placeholder rows above the data (old rows 12:13) ---
$ws.Rows("12:13").Delete()

# --- Drop the now-unused trailing columns E:F ---
$ws.Columns("E:F").Delete()

# --- Resize remaining columns: A:B to 20.710625, C:D to 25.710625 ---
$ws.Columns("A").ColumnWidth = 19.8
$ws.Columns("B").ColumnWidth = 19.8
$ws.Columns("C").ColumnWidth = 24.8
$ws.Columns("D").ColumnWidth = 24.8

# --- Keep the sheet's used range anchored at A1 (matches original layout) ---
$ws.Range("A1").Font.Bold = $false
